$d = $word.ActiveDocument

# 1. "Java EE, padrão MVC" -> "Java, padrão MVC" (drop the italic " EE" run)
$d.Content.Find.Execute("Java EE, padrão MVC", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Java, padrão MVC", 2)

# 2. "novas idéias, nova funcionalidade" -> "novas ideias, nova funcionalidade"
$d.Content.Find.Execute("novas idéias, nova funcionalidade", $true, $false, $false, $false, $false, `
    $true, 1, $false, "novas ideias, nova funcionalidade", 2)

# 3. "documentação, pois sendo bem formulada será mais fácil o seu entendimento" ->
#    "documentação, pois sendo bem formulada serão mais fáceis o seu entendimento"
$d.Content.Find.Execute("documentação, pois sendo bem formulada será mais fácil o seu entendimento", $true, $false, $false, $false, $false, `
    $true, 1, $false, "documentação, pois sendo bem formulada serão mais fáceis o seu entendimento", 2)

# 4. ", este artigo estará apresentado " -> ", este artigo vai apresentar "
$d.Content.Find.Execute(", este artigo estará apresentado ", $true, $false, $false, $false, $false, `
    $true, 1, $false, ", este artigo vai apresentar ", 2)

# 5. "Porém, no seu inicio, não havia" -> "Porém, no seu início, não havia"
$d.Content.Find.Execute("Porém, no seu inicio, não havia", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Porém, no seu início, não havia", 2)

# 6. "uma pagina exclusiva no site do" -> "uma página exclusiva no site do"
$d.Content.Find.Execute("uma pagina exclusiva no site do", $true, $false, $false, $false, $false, `
    $true, 1, $false, "uma página exclusiva no site do", 2)

# 7. "qualquer tipo de  projeto" (double space) -> "qualquer tipo de projeto"
$d.Content.Find.Execute("qualquer tipo de  projeto", $true, $false, $false, $false, $false, `
    $true, 1, $false, "qualquer tipo de projeto", 2)

# 8. "versão 4 ,trazendo" -> "versão 4, trazendo"
$d.Content.Find.Execute("versão 4 ,trazendo", $true, $false, $false, $false, $false, `
    $true, 1, $false, "versão 4, trazendo", 2)

# 9. "freqüentemente usando APIs" -> "frequentemente usando APIs"
$d.Content.Find.Execute("freqüentemente usando APIs", $true, $false, $false, $false, $false, `
    $true, 1, $false, "frequentemente usando APIs", 2)

# 10. "visualização, e agregação" -> "visualização e agregação"
$d.Content.Find.Execute("visualização, e agregação", $true, $false, $false, $false, $false, `
    $true, 1, $false, "visualização e agregação", 2)

# 11. "pessoal e profissional.”" -> "pessoal e profissional. ”"
$d.Content.Find.Execute("pessoal e profissional.”", $true, $false, $false, $false, $false, `
    $true, 1, $false, "pessoal e profissional. ”", 2)

# 12. Italicize the word "mashup"
$r = $d.Content
$r.Find.Execute("mashup", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Italic = 1
}

# 13. "que suportam as requisições de entrada através de  controladores" (double space) ->
#     "que suportam as requisições de entrada através de controladores"
$d.Content.Find.Execute("através de  controladores", $true, $false, $false, $false, $false, `
    $true, 1, $false, "através de controladores", 2)

# 14. "mais específicas do projetos conseguiram" -> "mais específicas dos projetos conseguiram"
$d.Content.Find.Execute("mais específicas do projetos conseguiram", $true, $false, $false, $false, $false, `
    $true, 1, $false, "mais específicas dos projetos conseguiram", 2)
